# Shift the dates in column F (rows 2-7) forward by one day,
# preserving existing cell style/number formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2..7) {
    $cell = $ws.Cells.Item($row, 6)   # column F
    $cell.Value2 = $cell.Value2 + 1
}
